# Add option to switch between SAFE and WIFM 4.1 at runtime
# Update the convergence analysis and more

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the WIFM 4.1 convergence-analysis inputs (rows 17-27 block) ---
# These are raw input values; all dependent formula cells (E19:E24, H19:H25)
# recalculate automatically from them.
$ws.Range("B22").Value = 97.021032000907951
$ws.Range("B23").Value = 97.024502896184202
$ws.Range("B25").Value = 12.3241714936403

# --- Relabel the coefficient row and add two new derived rows ---
# New shared strings must be introduced in the same order as the target
# workbook (rHeadDiff=30, COEFF=31) so the sharedStrings table matches.
$ws.Range("G26").Value = "rHeadDiff"
$ws.Range("H26").Formula = "=B22-B23"

$ws.Range("G25").Value = "COEFF"

$ws.Range("G27").Value = "Qsint"
$ws.Range("H27").Formula = "=H26*H25"

# --- Update the active selection on the sheet ---
$ws.Range("H28").Select()
